# Fellowship Conference 2022 poster - add Figure 1C to the caption textbox.
#
# 1) "Figures 1A & 1B: ..." -> "Figures 1A, 1B, & 1C: ..."
# 2) Split the trailing sentence so "be constrained." becomes its own run
#    reading "be constrained (Figure 1A)."

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# Locate the caption textbox ("Figures 1A & 1B: Graphical representations...")
$shape = $null
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $candidate = $s.Shapes.Item($i)
    if ($candidate.HasTextFrame -and $candidate.TextFrame.HasText) {
        if ($candidate.TextFrame.TextRange.Text.StartsWith("Figures 1A")) {
            $shape = $candidate
            break
        }
    }
}

$tr = $shape.TextFrame.TextRange

# --- Edit 1: "Figures 1A & 1B" -> "Figures 1A, 1B, & 1C" ------------------
$oldLead = "Figures 1A & 1B: Graphical representations of listeners responses to the "
$newLead = "Figures 1A, 1B, & 1C: Graphical representations of listeners responses to the "
$leadRange = $tr.Characters(1, $oldLead.Length)
$leadRange.Text = $newLead

# --- Edit 2: split "...will be constrained." into its own trailing run ----
$fullText = $tr.Text
$oldTail = "be constrained."
$tailIdx = $fullText.IndexOf($oldTail)
$tailStart = $tailIdx + 1
$tailRange = $tr.Characters($tailStart, $oldTail.Length)
$tailRange.Text = "be constrained (Figure 1A)."
